# Daily auto push: insert the new 2026/01/27 23:00 reading as a new row
# right after the existing 2026/01/27 block (row 735), pushing every
# following row down by one. Mirrors the author's commit:
#   "daily auto push: 2026-01-27 18:55 UTC"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 736..777 down to 737..778, creating a blank row 736.
$ws.Rows.Item(736).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/12/29"), not
# real Excel date serials. Force the cell to Text format first so the
# "2026/01/27" string assignment isn't auto-converted into a date value,
# then restore the default "Normal" style so it matches the formatting
# of every other data cell in the column.
$ws.Range("A736").NumberFormat = "@"
$ws.Range("A736").Value = "2026/01/27"
$ws.Range("A736").Style = "Normal"

$ws.Range("B736").Value = "火"
$ws.Range("C736").Value = 23
$ws.Range("D736").Value = 201
